$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix every header/value cell with "%" (the pid column switches from a
# numeric type to a "%"-prefixed string, same as every other column).
$ws.Range("A1").Value = "%pid"
$ws.Range("B1").Value = "%name"
$ws.Range("C1").Value = "%email"

$ws.Range("A2").Value = "%1"
$ws.Range("B2").Value = "%Foo Bar"
$ws.Range("C2").Value = "%foo@bar.com"

$ws.Range("A3").Value = "%2"
$ws.Range("B3").Value = "%Baz Quux"
$ws.Range("C3").Value = "%baz@quux.com"

# Drop the hyperlink on C2 (foo@bar.com), keep the one on C3 (baz@quux.com).
# Hyperlink objects obtained through the collection's enumerator support
# .Delete() individually (unlike re-fetching via .Item(n)), so collect them
# first and delete just the one we no longer want.
$hyperlinks = @()
foreach ($h in $ws.Hyperlinks) { $hyperlinks += $h }
$hyperlinks[0].Delete()

# Move the saved selection from B4 to C4.
$ws.Range("C4").Select()
